$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error metrics (ifo GDP component analysis preprocessing)
$data = @{
    2  = @(-0.07937728594090611, 1.316068219616539, 9.638393705627079, 3.104576252184359, 3.13384069789644, 52)
    3  = @(-0.02448439541659127, 1.462623545219845, 9.051394970963811, 3.008553634383773, 3.038389605539716, 51)
    4  = @(-0.08143474300426477, 1.364967268314127, 6.831940587747885, 2.613798115338651, 2.639053048004792, 50)
    5  = @(-0.001877683171190082, 1.602456891704275, 9.4526324803627, 3.074513372936065, 3.106373886794302, 49)
    6  = @(-0.1039814330460372, 1.471754775803563, 8.504042060525727, 2.91616907269207, 2.945154824746354, 48)
    7  = @(-0.1015771582092904, 1.769024061707294, 10.29407639076258, 3.208438310262889, 3.252319194407661, 36)
    8  = @(-0.1097886217413366, 1.815971081873004, 10.93971098436304, 3.307523391355387, 3.353961707055734, 35)
    9  = @(0.007295727682830271, 2.609854038779015, 19.14435649662439, 4.375426435974485, 4.50227040007513, 18)
    10 = @(-0.6545140871090607, 2.980052884870046, 23.04619067002914, 4.800644818149864, 4.987943682298893, 11)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = $vals[0]
    $ws.Cells.Item($rowNum, 3).Value = $vals[1]
    $ws.Cells.Item($rowNum, 4).Value = $vals[2]
    $ws.Cells.Item($rowNum, 5).Value = $vals[3]
    $ws.Cells.Item($rowNum, 6).Value = $vals[4]
    $ws.Cells.Item($rowNum, 7).Value = $vals[5]
}
